$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at the top of the data (rows 2-7), pushing existing
# data down by 6 rows (old row 2 -> row 8, ..., old row 321 -> row 327).
[void]$ws.Rows("2:7").Insert()

# New rows of ism_prices_paid data (most recent months prepended to the
# existing series), matching column A's existing date format (m/d/yy).
$ws.Range("A2:A7").NumberFormat = "m/d/yy"

$ws.Range("A2").Value = 45535
$ws.Range("B2").Value = 57.3

$ws.Range("A3").Value = 45504
$ws.Range("B3").Value = 57

$ws.Range("A4").Value = 45473
$ws.Range("B4").Value = 56.3

$ws.Range("A5").Value = 45443
$ws.Range("B5").Value = 58.1

$ws.Range("A6").Value = 45412
$ws.Range("B6").Value = 59.2

$ws.Range("A7").Value = 45382
$ws.Range("B7").Value = 53.4

# Resize column A to fit the date values now that the sheet has new data.
[void]$ws.Columns("A").AutoFit()

# Move the active selection to B2 and scroll the view back to the top.
[void]$ws.Range("B2").Select()
